$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Five Brothers
$ws.Range("B2").Value = "RET-19352"
$ws.Range("C2").Value = "Five Brothers"
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = 1746161400

# Row 3 - Sohan Electric & Hardware
$ws.Range("A3").Value = "Mugdho Corporation"
$ws.Range("B3").Value = "RET-35280"
$ws.Range("C3").Value = "Sohan Electric & Hardware"
$ws.Range("E3").Value = "Bkash "
$ws.Range("F3").Value = 1860207883

# Row 4 - Azim Mobile Center
$ws.Range("A4").Value = "Mugdho Corporation"
$ws.Range("B4").Value = "RET-36165"
$ws.Range("C4").Value = "Azim Mobile Center"
$ws.Range("E4").Value = "Bkash "
$ws.Range("F4").Value = 1764994148

# Row 5 - Gourango Hardware
$ws.Range("A5").Value = "Mugdho Corporation"
$ws.Range("B5").Value = "RET-36167"
$ws.Range("C5").Value = "Gourango Hardware"
$ws.Range("E5").Value = "Bkash "
$ws.Range("F5").Value = 1722309632
$ws.Range("G5").Clear()

# Row 6 - Bishakhi Enterprise (string order: C before B, per original authoring)
$ws.Range("A6").Value = "Mugdho Corporation"
$ws.Range("C6").Value = "Bishakhi Enterprise"
$ws.Range("B6").Value = "RET-34130"
$ws.Range("E6").Value = "Bkash "
$ws.Range("F6").Value = 1718898690

# Row 7 - CD Sound & Electronics (string order: C before B, per original authoring)
$ws.Range("A7").Value = "Mugdho Corporation"
$ws.Range("C7").Value = "CD Sound & Electronics"
$ws.Range("B7").Value = "RET-34136"
$ws.Range("E7").Value = "Bkash "
$ws.Range("F7").Value = 1730430130

# L12 gets a lone space value
$ws.Range("L12").Value = " "

# Column D is now hidden instead of auto-fit
$ws.Columns("D").Hidden = $true
